$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 89707
$ws.Range("B2").Value = "Sr. Matheus Santos"
$ws.Range("C2").Value = "Jurídico"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45101
$ws.Range("G2").Value = 8361.719999999999

# Row 3
$ws.Range("A3").Value = 86022
$ws.Range("B3").Value = "Carolina Pereira"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45105
$ws.Range("G3").Value = 12320.24

# Row 4
$ws.Range("A4").Value = 86504
$ws.Range("B4").Value = "Pedro Lucas Jesus"
$ws.Range("C4").Value = "Operações"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45083
$ws.Range("G4").Value = 4238.86

# Row 5
$ws.Range("A5").Value = 88438
$ws.Range("B5").Value = "Sra. Mirella Nascimento"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45094
$ws.Range("G5").Value = 3330.22

# Row 6
$ws.Range("A6").Value = 21566
$ws.Range("B6").Value = "Srta. Laís Pereira"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 4250.35

# Row 7
$ws.Range("A7").Value = 94764
$ws.Range("B7").Value = "Srta. Bianca Almeida"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45095
$ws.Range("G7").Value = 3651.46

# Row 8
$ws.Range("A8").Value = 55288
$ws.Range("B8").Value = "Srta. Ana Sophia Silveira"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45095
$ws.Range("G8").Value = 8327.99

# Row 9
$ws.Range("A9").Value = 65695
$ws.Range("B9").Value = "Isabel Souza"
$ws.Range("C9").Value = "Operações"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45097
$ws.Range("G9").Value = 10803.59

# Row 10
$ws.Range("A10").Value = 91290
$ws.Range("B10").Value = "Thomas Moura"
$ws.Range("C10").Value = "TI"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 45098
$ws.Range("G10").Value = 3597.36

# Row 11
$ws.Range("A11").Value = 80468
$ws.Range("B11").Value = "Milena Pereira"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45080
$ws.Range("G11").Value = 10682.54
